$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 4.75
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = 2.6
$ws.Range("L2").Value = 5.5
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.63
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 9.5
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 13
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 10
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 41
$ws.Range("AN2").Value = 3.75
$ws.Range("AQ2").Value = 19
$ws.Range("AR2").Value = 34
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 3.75
$ws.Range("AU2").Value = 7.5
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 401
$ws.Range("AX2").Value = 7.5
$ws.Range("AZ2").Value = 29
$ws.Range("BA2").Value = 81
$ws.Range("BB2").Value = 81
$ws.Range("BC2").Value = 151
# Row 4
$ws.Range("G4").Value = 2.1
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.88
$ws.Range("L4").Value = 4.33
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.8
$ws.Range("X4").Value = 9.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 19
$ws.Range("AC4").Value = 8
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 9.5
$ws.Range("AI4").Value = 17
$ws.Range("AJ4").Value = 13
$ws.Range("AK4").Value = 41
$ws.Range("AN4").Value = 4
$ws.Range("AO4").Value = 12
$ws.Range("AT4").Value = 2.5
$ws.Range("AV4").Value = 67
$ws.Range("AX4").Value = 5.5
$ws.Range("AY4").Value = 21
$ws.Range("BB4").Value = 101
# Row 5
$ws.Range("G5").Value = 32
$ws.Range("H5").Value = 9.5
$ws.Range("J5").Value = 21
$ws.Range("K5").Value = 4.15
$ws.Range("L5").Value = 1.19
$ws.Range("S5").Value = 1.05
$ws.Range("T5").Value = 7.8
$ws.Range("U5").Value = 1.9
$ws.Range("V5").Value = 1.86
$ws.Range("W5").Value = 300
$ws.Range("Y5").Value = 200
$ws.Range("AB5").Value = 350
$ws.Range("AC5").Value = 45
$ws.Range("AD5").Value = 37
$ws.Range("AH5").Value = 23
$ws.Range("AI5").Value = 11.25
$ws.Range("AJ5").Value = 17.5
$ws.Range("AK5").Value = 7.9
$ws.Range("AM5").Value = 32
$ws.Range("AN5").Value = 40
$ws.Range("AO5").Value = 350
$ws.Range("AP5").Value = 120
$ws.Range("AT5").Value = 7.4
$ws.Range("AU5").Value = 13
$ws.Range("AV5").Value = 70
$ws.Range("AX5").Value = 4.05
$ws.Range("AY5").Value = 4
$ws.Range("AZ5").Value = 12.5
$ws.Range("BB5").Value = 17.5
$ws.Range("BC5").Value = 100
# Row 6
$ws.Range("H6").Value = 2.8
$ws.Range("I6").Value = 3.8
# Row 7
$ws.Range("M7").Value = 1.04
$ws.Range("O7").Value = 1.3
$ws.Range("AW7").Value = 126
# Row 8
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 1.47
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.48
# Row 9
$ws.Range("M9").Value = 1.11
$ws.Range("O9").Value = 1.63
# Row 15
$ws.Range("G15").Value = 1.48
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 2.05
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("W15").Value = 7
$ws.Range("AC15").Value = 11
$ws.Range("AH15").Value = 17
$ws.Range("AI15").Value = 34
$ws.Range("AO15").Value = 7.5
$ws.Range("AV15").Value = 51
$ws.Range("BA15").Value = 126
# Row 21
$ws.Range("G21").Value = 2.05
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 3.75
$ws.Range("J21").Value = 2.88
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 4.33
$ws.Range("N21").Value = 7.5
$ws.Range("X21").Value = 9.5
$ws.Range("Z21").Value = 19
$ws.Range("AA21").Value = 19
$ws.Range("AE21").Value = 15
$ws.Range("AH21").Value = 9.5
$ws.Range("AJ21").Value = 13
$ws.Range("AO21").Value = 12
$ws.Range("AY21").Value = 21
$ws.Range("BA21").Value = 67
$ws.Range("BC21").Value = 251
# Row 31
$ws.Range("G31").Value = 2.92
$ws.Range("I31").Value = 1.95
$ws.Range("J31").Value = 3.25
$ws.Range("K31").Value = 2.45
$ws.Range("L31").Value = 2.4
$ws.Range("M31").Value = 1.02
$ws.Range("N31").Value = 13
$ws.Range("P31").Value = 4.9
$ws.Range("U31").Value = 1.41
$ws.Range("V31").Value = 2.81
$ws.Range("W31").Value = 12.5
$ws.Range("X31").Value = 16.5
$ws.Range("Y31").Value = 9.5
$ws.Range("Z31").Value = 30
$ws.Range("AA31").Value = 17.5
$ws.Range("AB31").Value = 18.5
$ws.Range("AC31").Value = 18.5
$ws.Range("AD31").Value = 7.4
$ws.Range("AE31").Value = 10
$ws.Range("AF31").Value = 27
$ws.Range("AG31").Value = 120
$ws.Range("AH31").Value = 10.25
$ws.Range("AI31").Value = 10.75
$ws.Range("AJ31").Value = 7.7
$ws.Range("AK31").Value = 16
$ws.Range("AL31").Value = 11.5
$ws.Range("AN31").Value = 5.5
$ws.Range("AO31").Value = 14.5
$ws.Range("AP31").Value = 17
$ws.Range("AQ31").Value = 55
$ws.Range("AR31").Value = 65
$ws.Range("AS31").Value = 150
$ws.Range("AT31").Value = 3.85
$ws.Range("AU31").Value = 6.3
$ws.Range("AV31").Value = 37
$ws.Range("AW31").Value = 450
$ws.Range("AX31").Value = 4.4
$ws.Range("AY31").Value = 9.5
$ws.Range("BA31").Value = 29
$ws.Range("BC31").Value = 110
# Row 46
$ws.Range("G46").Value = 2.32
$ws.Range("J46").Value = 2.95
$ws.Range("K46").Value = 2.1
$ws.Range("L46").Value = 3.3
$ws.Range("N46").Value = 11.5
$ws.Range("O46").Value = 1.21
$ws.Range("P46").Value = 3.55
$ws.Range("Q46").Value = 1.62
$ws.Range("R46").Value = 2.02
$ws.Range("U46").Value = 1.52
$ws.Range("V46").Value = 2.22
$ws.Range("X46").Value = 13
$ws.Range("Z46").Value = 25
$ws.Range("AB46").Value = 23
$ws.Range("AC46").Value = 12.5
$ws.Range("AD46").Value = 6.7
$ws.Range("AE46").Value = 11.5
$ws.Range("AF46").Value = 40
$ws.Range("AG46").Value = 250
$ws.Range("AH46").Value = 11.75
$ws.Range("AI46").Value = 17
$ws.Range("AL46").Value = 21
$ws.Range("AM46").Value = 24
$ws.Range("AN46").Value = 4.35
$ws.Range("AO46").Value = 12.5
$ws.Range("AS46").Value = 200
$ws.Range("AT46").Value = 2.62
$ws.Range("AU46").Value = 6.5
$ws.Range("AV46").Value = 50
$ws.Range("AX46").Value = 4.8
$ws.Range("AY46").Value = 14.5
$ws.Range("AZ46").Value = 19.5
$ws.Range("BA46").Value = 65
$ws.Range("BB46").Value = 80
$ws.Range("BC46").Value = 200
# Row 47
$ws.Range("I47").Value = 8.75
$ws.Range("J47").Value = 1.72
$ws.Range("L47").Value = 7.5
$ws.Range("R47").Value = 2.1
$ws.Range("U47").Value = 2.02
$ws.Range("V47").Value = 1.62
$ws.Range("AH47").Value = 23
$ws.Range("AL47").Value = 110
$ws.Range("AP47").Value = 16.5
$ws.Range("AR47").Value = 45
$ws.Range("AT47").Value = 3.1
$ws.Range("AU47").Value = 9.25
$ws.Range("AX47").Value = 9.25
# Row 61
$ws.Range("H61").Value = 3.25
$ws.Range("I61").Value = 4.1
$ws.Range("K61").Value = 2
$ws.Range("N61").Value = 7.5
$ws.Range("W61").Value = 5.5
$ws.Range("X61").Value = 7.5
$ws.Range("AC61").Value = 7.5
$ws.Range("AE61").Value = 21
$ws.Range("AK61").Value = 51
$ws.Range("AM61").Value = 51
$ws.Range("AN61").Value = 3.6
$ws.Range("AP61").Value = 26
$ws.Range("AZ61").Value = 41
$ws.Range("BB61").Value = 151
# Row 62
$ws.Range("K62").Value = 1.92
# Row 63
$ws.Range("M63").Value = 1.06
$ws.Range("N63").Value = 10
